$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. week49: fill in the Friday (row 8) begin/end times that were
#    previously blank, and move the selection to D9.
# ---------------------------------------------------------------
$week49 = $wb.Worksheets.Item("week49")
$week49.Range("C8").Value = 0.40625
$week49.Range("D8").Value = 0.5625

# ---------------------------------------------------------------
# 2. Add a new "week50" sheet, positioned right after "week49".
#    Copying week49 keeps all number formats / column widths /
#    row heights / styles identical to the other week sheets.
# ---------------------------------------------------------------
$week49.Copy($null, $week49)
$week50 = $wb.Worksheets.Item("week49 (2)")
$week50.Name = "week50"

# ---------------------------------------------------------------
# 3. Build the new "totaal" sheet from a copy of the existing one
#    (placed right after week50), then retire the old "totaal"
#    sheet (renamed to keep formatting continuity for week50 ids).
# ---------------------------------------------------------------
$totaalOld = $wb.Worksheets.Item("totaal")
$totaalOld.Copy($null, $week50)
$totaalNew = $wb.Worksheets.Item("totaal (2)")
$totaalOld.Delete()
$totaalNew.Name = "totaal"

# ---------------------------------------------------------------
# 4. Populate week50's activity rows with the new log entries.
# ---------------------------------------------------------------
$week50.Range("C7").Value = 0.364583333333333
$week50.Range("D7").Value = 0.370833333333333
$week50.Range("F7").Value = "Een nieuw tabblad toegevoegd aan het logboek en"

$week50.Range("C8").Value = 0.371527777777778
$week50.Range("C8").NumberFormat = "HH:MM:SS"
$week50.Range("D8").Value = 0.427083333333333
$week50.Range("F8").Value = "Sessionclass toegevoegd loginclass bijgewerkt"

# week50 becomes the active/visible sheet.
$week50.Activate()
$week50.Range("F8").Select()

# ---------------------------------------------------------------
# 5. Extend the "totaal" sheet with rows for week49 and week50,
#    and update the grand-total SUM range.
# ---------------------------------------------------------------
$totaal = $wb.Worksheets.Item("totaal")
$totaal.Rows.Item(9).Insert()
$totaal.Rows.Item(9).Insert()

$totaal.Range("A9").Value = 49
$totaal.Range("B9").Formula = "=week49!G19"

$totaal.Range("A10").Value = 50
$totaal.Range("B10").Formula = "=week50!G19"

$totaal.Range("B11").Formula = "=SUM(B7:B10)"
$totaal.Range("B12").Select()

# ---------------------------------------------------------------
# 6. week49 is no longer the active tab; move its selection too.
# ---------------------------------------------------------------
$week49.Range("D9").Select()

# Make sure week50 stays the tab that is shown on open.
$week50.Activate()
